$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift timestamps: the newest block (rows 2-15) gets a brand new
# timestamp, and the older blocks inherit the previous block's value
# (rows 16-29 take the old rows 2-15 value, rows 30-43 take the old
# rows 16-29 value), matching an automated "Actualizar" refresh run.

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = 44266.54345219229
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = 44266.52208146991
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = 44266.50069510417
}
